# Updated cryptos list - apply price (D) and volume% (E) changes per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "26.721.49"
$ws.Range("E2").Value = "  +4.42%  "

Set-TextValue $ws.Range("D3") "1.874.81"
$ws.Range("E3").Value = "  +3.62%  "

$ws.Range("E4").Value = "  -0.10%  "

Set-TextValue $ws.Range("D5") "282.63"
$ws.Range("E5").Value = "  +2.01%  "

Set-TextValue $ws.Range("D6") "0.9999"
$ws.Range("E6").Value = "  -0.15%  "

Set-TextValue $ws.Range("D7") "0.5176"
$ws.Range("E7").Value = "  +3.12%  "

Set-TextValue $ws.Range("D8") "0.3531"
$ws.Range("E8").Value = "  +0.85%  "

Set-TextValue $ws.Range("D9") "45.29"
$ws.Range("E9").Value = "  +3.28%  "

Set-TextValue $ws.Range("D10") "0.07120"
$ws.Range("E10").Value = "  +7.03%  "

Set-TextValue $ws.Range("D11") "20.23"
$ws.Range("E11").Value = "  +1.00%  "

Set-TextValue $ws.Range("D12") "0.8207"
$ws.Range("E12").Value = "  -1.89%  "

Set-TextValue $ws.Range("D13") "0.07764"
$ws.Range("E13").Value = "  -0.76%  "

Set-TextValue $ws.Range("D14") "1.866.42"
$ws.Range("E14").Value = "  +3.23%  "

Set-TextValue $ws.Range("D15") "5.173"
$ws.Range("E15").Value = "  +2.70%  "

Set-TextValue $ws.Range("D16") "89.73"
$ws.Range("E16").Value = "  +2.86%  "

Set-TextValue $ws.Range("D17") "1.000"
$ws.Range("E17").Value = "  -0.03%  "

Set-TextValue $ws.Range("D18") "14.45"
$ws.Range("E18").Value = "  +4.05%  "

Set-TextValue $ws.Range("D19") "0.000008172"
$ws.Range("E19").Value = "  +3.78%  "

$ws.Range("E20").Value = "  -0.18%  "

Set-TextValue $ws.Range("D21") "26.756.79"
$ws.Range("E21").Value = "  +4.31%  "

Set-TextValue $ws.Range("D22") "4.801"
$ws.Range("E22").Value = "  +1.84%  "

Set-TextValue $ws.Range("D23") "10.18"
$ws.Range("E23").Value = "  +2.22%  "

Set-TextValue $ws.Range("D24") "6.248"
$ws.Range("E24").Value = "  +3.22%  "

Set-TextValue $ws.Range("D25") "2.435"
$ws.Range("E25").Value = "  +15.51%  "

Set-TextValue $ws.Range("D26") "145.78"
$ws.Range("E26").Value = "  +3.37%  "

$ws.Range("E27").Value = "  +3.31%  "

Set-TextValue $ws.Range("D28") "1.663"
$ws.Range("E28").Value = "  -0.11%  "

Set-TextValue $ws.Range("D29") "111.40"
$ws.Range("E29").Value = "  +2.49%  "

Set-TextValue $ws.Range("D30") "4.419"
$ws.Range("E30").Value = "  +2.98%  "

Set-TextValue $ws.Range("D31") "4.365"
$ws.Range("E31").Value = "  +3.83%  "

Set-TextValue $ws.Range("D32") "0.08859"
$ws.Range("E32").Value = "  +0.47%  "

Set-TextValue $ws.Range("D33") "0.04917"
$ws.Range("E33").Value = "  +2.63%  "

$ws.Range("E34").Value = "  +5.14%  "

$ws.Range("E35").Value = "  +2.11%  "

Set-TextValue $ws.Range("D36") "3.300"
$ws.Range("E36").Value = "  +8.72%  "

Set-TextValue $ws.Range("D37") "2.867"
$ws.Range("E37").Value = "  -0.03%  "

Set-TextValue $ws.Range("D38") "2.426"
$ws.Range("E38").Value = "  +7.49%  "

Set-TextValue $ws.Range("D39") "0.5316"
$ws.Range("E39").Value = "  +2.08%  "

Set-TextValue $ws.Range("D40") "0.01889"
$ws.Range("E40").Value = "  +1.63%  "

Set-TextValue $ws.Range("D41") "0.9759"
$ws.Range("E41").Value = "  +1.85%  "

Set-TextValue $ws.Range("D42") "116.59"
$ws.Range("E42").Value = "  +4.55%  "

Set-TextValue $ws.Range("D43") "6.312"
$ws.Range("E43").Value = "  +2.27%  "

Set-TextValue $ws.Range("D44") "8.196"
$ws.Range("E44").Value = "  +1.87%  "

Set-TextValue $ws.Range("D45") "0.9995"
$ws.Range("E45").Value = "  -0.15%  "

Set-TextValue $ws.Range("D46") "0.4621"
$ws.Range("E46").Value = "  +0.59%  "

Set-TextValue $ws.Range("D47") "0.1369"
$ws.Range("E47").Value = "  -0.65%  "

Set-TextValue $ws.Range("D48") "9.541"
$ws.Range("E48").Value = "  +3.93%  "

Set-TextValue $ws.Range("D49") "36.70"
$ws.Range("E49").Value = "  +3.17%  "

Set-TextValue $ws.Range("D50") "1.518"
$ws.Range("E50").Value = "  +2.15%  "

Set-TextValue $ws.Range("D51") "0.05935"
$ws.Range("E51").Value = "  +1.52%  "
